$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.266.53'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '3.787.89'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '432.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.52%  '
$ws.Range('E7').Value = '  +1.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000311'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -14.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.58'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.40'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').Value = '4.391.27'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.94'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '3.766.05'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('E19').Value = '  +5.22%  '
$ws.Range('D20').Value = '66.317.54'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '406.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('E23').Value = '  +6.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '36.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +31.70%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.10%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.87'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.138'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +12.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '704.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.96%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.148'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.09%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +28.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0472'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +37.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.141'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0672'
$ws.Range('E44').Value = '  -9.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.328'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '140.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.72%  '
